# "Received Qorvo LNA boards" - update the Boards inventory sheet to mark
# the three Qorvo TQP3M90xx LNA boards as received (3 each) instead of
# still on order.

$wb = $excel.ActiveWorkbook
$wsBoards = $wb.Worksheets.Item("Boards")
$wsInventory = $wb.Worksheets.Item("Inventory")

$receivedRows = 6, 7, 8
foreach ($r in $receivedRows) {
    $statusCell = $wsBoards.Cells.Item($r, 2)
    $statusCell.Value = "RCVD"
    $statusCell.Interior.Color = 5296274
    $statusCell.HorizontalAlignment = -4108

    $wsBoards.Cells.Item($r, 3).Value = 3
}

# Restore the cursor positions Excel leaves behind after the edit.
[void]$wsBoards.Range("A17").Select()
[void]$wsInventory.Range("D2").Select()
